$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.256.29"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3
$ws.Range("D3").Value = "1.790.36"
$ws.Range("E3").Value = "  +0.04%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'226.05"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6
$ws.Range("E6").Value = "  +0.34%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'32.31"
$ws.Range("E8").Value = "  +0.50%  "

# Row 9
$ws.Range("E9").Value = "  +0.28%  "

# Row 10
$ws.Range("D10").Value = "'0.0691"
$ws.Range("E10").Value = "  -0.34%  "

# Row 11
$ws.Range("D11").Value = "'0.0947"
$ws.Range("E11").Value = "  +0.70%  "

# Row 12
$ws.Range("D12").Value = "2.049.20"
$ws.Range("E12").Value = "  +0.12%  "

# Row 13
$ws.Range("E13").Value = "  -3.22%  "

# Row 14
$ws.Range("D14").Value = "1.793.57"
$ws.Range("E14").Value = "  +0.16%  "

# Row 15
$ws.Range("E15").Value = "  +0.50%  "

# Row 16
$ws.Range("D16").Value = "34.240.93"
$ws.Range("E16").Value = "  +0.45%  "

# Row 17
$ws.Range("E17").Value = "  +0.21%  "

# Row 18
$ws.Range("D18").Value = "'68.05"
$ws.Range("E18").Value = "  +0.03%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0807"
$ws.Range("E19").Value = "  +3.13%  "

# Row 20
$ws.Range("D20").Value = "'246.52"
$ws.Range("E20").Value = "  +0.75%  "

# Row 21
$ws.Range("D21").Value = "'10.99"
$ws.Range("E21").Value = "  +0.57%  "

# Row 22
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
$ws.Range("D23").Value = "'4.20"
$ws.Range("E23").Value = "  +2.30%  "

# Row 24
$ws.Range("D24").Value = "'2.06"
$ws.Range("E24").Value = "  +0.66%  "

# Row 25
$ws.Range("D25").Value = "'161.54"
$ws.Range("E25").Value = "  -0.59%  "

# Row 26
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("D27").Value = "'16.34"
$ws.Range("E27").Value = "  +0.29%  "

# Row 28
$ws.Range("E28").Value = "  +0.79%  "

# Row 29
$ws.Range("E29").Value = "  +0.19%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.0521"
$ws.Range("E30").Value = "  +0.06%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.23"
$ws.Range("E31").Value = "  -0.31%  "

# Row 32
$ws.Range("D32").Value = "'3.75"
$ws.Range("E32").Value = "  +2.37%  "

# Row 33
$ws.Range("D33").Value = "'3.75"
$ws.Range("E33").Value = "  +3.82%  "

# Row 34
$ws.Range("E34").Value = "  -1.41%  "

# Row 35
$ws.Range("D35").Value = "1.444.68"
$ws.Range("E35").Value = "  +1.83%  "

# Row 36
$ws.Range("D36").Value = "'2.59"
$ws.Range("E36").Value = "  +10.33%  "

# Row 37
$ws.Range("D37").Value = "'0.669"
$ws.Range("E37").Value = "  +4.13%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.05"
$ws.Range("E38").Value = "  +1.13%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0191"
$ws.Range("E39").Value = "  -0.19%  "

# Row 40
$ws.Range("D40").Value = "'82.05"
$ws.Range("E40").Value = "  +1.55%  "

# Row 41
$ws.Range("D41").Value = "'14.10"
$ws.Range("E41").Value = "  +5.03%  "

# Row 42
$ws.Range("D42").Value = "'2.39"
$ws.Range("E42").Value = "  +1.17%  "

# Row 43
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.73"
$ws.Range("E43").Value = "  +1.47%  "

# Row 44
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'0.923"
$ws.Range("E44").Value = "  +0.56%  "

# Row 45
$ws.Range("E45").Value = "  +2.08%  "

# Row 46
$ws.Range("D46").Value = "'6.08"
$ws.Range("E46").Value = "  +0.24%  "

# Row 47
$ws.Range("E47").Value = "  +0.30%  "

# Row 48
$ws.Range("D48").Value = "1.948.74"
$ws.Range("E48").Value = "  +0.04%  "

# Row 49
$ws.Range("D49").Value = "'105.50"
$ws.Range("E49").Value = "  -1.77%  "

# Row 50
$ws.Range("E50").Value = "  -0.03%  "

# Row 51
$ws.Range("D51").Value = "0.0₆0129"
$ws.Range("E51").Value = "  -5.94%  "
